$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 48.4
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 91.8
$ws.Range("K6").Value = 15
$ws.Range("L6").Value = 275.4
$ws.Range("M6").Value = 97
$ws.Range("N6").Value = -499.4

$ws.Range("H69").Value = 7249.6855
$ws.Range("I69").Value = 2999.5
$ws.Range("J69").Value = 7507.273
$ws.Range("K69").Value = 8998.5
$ws.Range("L69").Value = 22521.819
$ws.Range("M69").Value = -8124.5
$ws.Range("N69").Value = -24269.819

$ws.Range("H72").Value = 7249.6855
$ws.Range("I72").Value = 2999.5
$ws.Range("J72").Value = 7507.273
$ws.Range("K72").Value = 26995.5
$ws.Range("L72").Value = 67565.45699999999
$ws.Range("M72").Value = -22627.5
$ws.Range("N72").Value = -76301.45699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4401.4736
$ws.Range("I74").Value = 4374.6665
$ws.Range("J74").Value = 4502
$ws.Range("K74").Value = 4374.6665
$ws.Range("L74").Value = 4502
$ws.Range("M74").Value = -3500.6665
$ws.Range("N74").Value = -6250

$ws.Range("H77").Value = 4401.4736
$ws.Range("I77").Value = 4374.6665
$ws.Range("J77").Value = 4502
$ws.Range("K77").Value = 21873.3325
$ws.Range("L77").Value = 22510
$ws.Range("M77").Value = -17505.3325
$ws.Range("N77").Value = -31246

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H97").Value = 986.6667
$ws.Range("I97").Value = 723.1111
$ws.Range("K97").Value = 723.1111
$ws.Range("M97").Value = -227.1111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4203.1665
$ws.Range("I20").Value = 7333.3335
$ws.Range("K20").Value = 7333.3335
$ws.Range("M20").Value = -7086.3335

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H134").Value = 1470.3
$ws.Range("I134").Value = 1247.7368
$ws.Range("K134").Value = 3743.2104
$ws.Range("M134").Value = -1208.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 970
$ws.Range("I16").Value = 970
$ws.Range("K16").Value = 970
$ws.Range("M16").Value = -683

$ws.Range("H22").Value = 2572.4285
$ws.Range("I22").Value = 1251.75
$ws.Range("J22").Value = 4333.3335
$ws.Range("K22").Value = 1251.75
$ws.Range("L22").Value = 4333.3335
$ws.Range("M22").Value = -901.75
$ws.Range("N22").Value = -5033.3335

$ws.Range("H58").Value = 2017.8
$ws.Range("I58").Value = 1897.25
$ws.Range("K58").Value = 1897.25
$ws.Range("M58").Value = -1694.25

$ws.Range("H99").Value = 2523.6
$ws.Range("I99").Value = 1702.6666
$ws.Range("K99").Value = 1702.6666
$ws.Range("M99").Value = -204.6666

$ws.Range("H113").Value = 970
$ws.Range("I113").Value = 970
$ws.Range("K113").Value = 970
$ws.Range("M113").Value = 1200

$ws.Range("H126").Value = 2523.6
$ws.Range("I126").Value = 1702.6666
$ws.Range("K126").Value = 5107.9998
$ws.Range("M126").Value = -2637.9998

$ws.Range("H132").Value = 3274.2104
$ws.Range("I132").Value = 3167.2222
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 9501.6666
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -6971.6666
$ws.Range("N132").Value = -20660

$ws.Range("H136").Value = 2017.8
$ws.Range("I136").Value = 1897.25
$ws.Range("K136").Value = 5691.75
$ws.Range("M136").Value = -3141.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 779
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 779
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H107").Value = 466.7
$ws.Range("I107").Value = 418.66666
$ws.Range("J107").Value = 475.17648
$ws.Range("K107").Value = 1255.99998
$ws.Range("L107").Value = 1425.52944
$ws.Range("M107").Value = 664.0000199999999
$ws.Range("N107").Value = -5265.52944

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 8487.25
$ws.Range("J46").Value = 9474.5
$ws.Range("L46").Value = 9474.5
$ws.Range("N46").Value = -9786.5

$ws.Range("H97").Value = 901.1111
$ws.Range("I97").Value = 1117.8
$ws.Range("K97").Value = 1117.8
$ws.Range("M97").Value = -621.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 2997.5
$ws.Range("J17").Value = 2995
$ws.Range("L17").Value = 2995
$ws.Range("N17").Value = -3335

$ws.Range("H22").Value = 910.55554
$ws.Range("I22").Value = 742.1429000000001
$ws.Range("K22").Value = 742.1429000000001
$ws.Range("M22").Value = -447.1429000000001

$ws.Range("H27").Value = 910.55554
$ws.Range("I27").Value = 742.1429000000001
$ws.Range("K27").Value = 742.1429000000001
$ws.Range("M27").Value = -635.1429000000001

$ws.Range("H46").Value = 3071.8125
$ws.Range("I46").Value = 670.2
$ws.Range("K46").Value = 670.2
$ws.Range("M46").Value = -482.2

$ws.Range("H55").Value = 884.0952
$ws.Range("I55").Value = 989.7857
$ws.Range("J55").Value = 672.7143
$ws.Range("K55").Value = 989.7857
$ws.Range("L55").Value = 672.7143
$ws.Range("M55").Value = -816.7857
$ws.Range("N55").Value = -1018.7143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 54999.5
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

$ws.Range("H67").Value = 54999.5
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

$ws.Range("H132").Value = 1368.6
$ws.Range("I132").Value = 1368.6
$ws.Range("K132").Value = 4105.799999999999
$ws.Range("M132").Value = -1575.799999999999

$ws.Range("H136").Value = 2857.1667
$ws.Range("I136").Value = 2503.4285
$ws.Range("J136").Value = 5333.3335
$ws.Range("K136").Value = 7510.2855
$ws.Range("L136").Value = 16000.0005
$ws.Range("M136").Value = -4960.2855
$ws.Range("N136").Value = -21100.0005
